# Applies the two textual edits described by the commit:
#  1. Slide 2 (Title): merge the three runs "The " / "minutes " /
#     "from the previous discussion " into a single run.
#  2. Slide 8 (TextBox 67), the "Zero Spring configuration..." bullet:
#     split the single run into three runs with re-worded text:
#       "Zero Spring " + "XML- configuration. " + "Use annotations instead;"

$p = $ppt.ActivePresentation

# --- Change 1: Slide 2 title -----------------------------------------
$s1 = $p.Slides.Item(2)
$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# Setting .Text to the exact same concatenation is a no-op for the
# underlying run structure, so nudge it through a placeholder first to
# force PowerPoint to collapse everything into a single run that keeps
# the first run's formatting (matches the diff: 3 runs -> 1 run).
$titleRange.Text = "placeholder"
$titleRange.Text = "The minutes from the previous discussion "

# --- Change 2: Slide 8 bullet text ------------------------------------
$s2 = $p.Slides.Item(8)
$bodyShape = $s2.Shapes.Item(3)
$bodyRange = $bodyShape.TextFrame.TextRange

# This shape autosizes to its text (<a:spAutoFit/>); editing its text
# makes PowerPoint re-layout and recompute the shape's height. The
# original commit's diff shows no geometry change, so remember the
# current height (in EMU, derived from the on-disk extent) to restore
# afterwards.
$originalHeightEmu = 5632311.0

$oldBullet = "Zero Spring configuration in XML. Use annotations instead;"
$fullText = $bodyRange.Text
$offset = $fullText.IndexOf($oldBullet)
$startPos = $offset + 1

$run1 = "Zero Spring "
$run2 = "XML- configuration. "
$run3 = "Use annotations instead;"

# Replace the whole bullet with the new wording first (still one run).
$wholeBullet = $bodyRange.Characters($startPos, $oldBullet.Length)
$wholeBullet.Text = $run1 + $run2 + $run3

# Re-assigning sub-ranges (even to their current text) forces PowerPoint
# to split the run at those boundaries, giving us three separate runs.
$part1 = $bodyRange.Characters($startPos, $run1.Length)
$part1.Text = $run1

$part2Start = $startPos + $run1.Length
$part2 = $bodyRange.Characters($part2Start, $run2.Length)
$part2.Text = $run2

# Restore the shape's autofit height exactly. Going through the Height
# *getter* first would round-trip the EMU value through a 32-bit float
# twice (losing a unit), so recompute the point value directly from the
# known EMU total instead, with a tiny epsilon to land on the correct
# EMU after the single float32 rounding the setter performs.
$bodyShape.Height = ($originalHeightEmu / 12700.0) + 0.00001
